{"js": "// Replace each two-digit multiplication equation in the document's table\n// cells with its new value, as described by the diff. Every old value is\n// unique in the document, so a direct search + replace per pair is safe.\nconst replacements = [\n  [\"98\u00d793=9114\", \"71\u00d729=2059\"],\n  [\"21\u00d714=294\", \"59\u00d724=1416\"],\n  [\"36\u00d714=504\", \"33\u00d792=3036\"],\n  [\"39\u00d743=1677\", \"41\u00d721=861\"],\n  [\"82\u00d718=1476\", \"52\u00d759=3068\"],\n  [\"13\u00d792=1196\", \"41\u00d716=656\"],\n  [\"23\u00d779=1817\", \"91\u00d718=1638\"],\n  [\"59\u00d775=4425\", \"61\u00d772=4392\"],\n  [\"85\u00d733=2805\", \"50\u00d765=3250\"],\n  [\"15\u00d757=855\", \"96\u00d733=3168\"],\n  [\"39\u00d734=1326\", \"55\u00d725=1375\"],\n  [\"94\u00d720=1880\", \"93\u00d736=3348\"],\n  [\"78\u00d736=2808\", \"54\u00d766=3564\"],\n  [\"45\u00d734=1530\", \"98\u00d711=1078\"],\n  [\"25\u00d756=1400\", \"18\u00d752=936\"],\n  [\"76\u00d711=836\", \"77\u00d796=7392\"],\n  [\"65\u00d711=715\", \"36\u00d763=2268\"],\n  [\"72\u00d787=6264\", \"26\u00d774=1924\"],\n  [\"15\u00d755=825\", \"17\u00d793=1581\"],\n  [\"30\u00d734=1020\", \"19\u00d772=1368\"],\n  [\"66\u00d761=4026\", \"30\u00d724=720\"],\n  [\"43\u00d744=1892\", \"97\u00d783=8051\"],\n  [\"16\u00d761=976\", \"40\u00d715=600\"],\n  [\"43\u00d794=4042\", \"41\u00d722=902\"],\n  [\"86\u00d783=7138\", \"33\u00d768=2244\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Pattern not found: \" + oldText);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication equation in the document's table\n# cells with its new value, as described by the diff. Every old value is\n# unique in the document, so a direct Find/Replace per pair is safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"98\u00d793=9114\", \"71\u00d729=2059\"),\n    @(\"21\u00d714=294\", \"59\u00d724=1416\"),\n    @(\"36\u00d714=504\", \"33\u00d792=3036\"),\n    @(\"39\u00d743=1677\", \"41\u00d721=861\"),\n    @(\"82\u00d718=1476\", \"52\u00d759=3068\"),\n    @(\"13\u00d792=1196\", \"41\u00d716=656\"),\n    @(\"23\u00d779=1817\", \"91\u00d718=1638\"),\n    @(\"59\u00d775=4425\", \"61\u00d772=4392\"),\n    @(\"85\u00d733=2805\", \"50\u00d765=3250\"),\n    @(\"15\u00d757=855\", \"96\u00d733=3168\"),\n    @(\"39\u00d734=1326\", \"55\u00d725=1375\"),\n    @(\"94\u00d720=1880\", \"93\u00d736=3348\"),\n    @(\"78\u00d736=2808\", \"54\u00d766=3564\"),\n    @(\"45\u00d734=1530\", \"98\u00d711=1078\"),\n    @(\"25\u00d756=1400\", \"18\u00d752=936\"),\n    @(\"76\u00d711=836\", \"77\u00d796=7392\"),\n    @(\"65\u00d711=715\", \"36\u00d763=2268\"),\n    @(\"72\u00d787=6264\", \"26\u00d774=1924\"),\n    @(\"15\u00d755=825\", \"17\u00d793=1581\"),\n    @(\"30\u00d734=1020\", \"19\u00d772=1368\"),\n    @(\"66\u00d761=4026\", \"30\u00d724=720\"),\n    @(\"43\u00d744=1892\", \"97\u00d783=8051\"),\n    @(\"16\u00d761=976\", \"40\u00d715=600\"),\n    @(\"43\u00d794=4042\", \"41\u00d722=902\"),\n    @(\"86\u00d783=7138\", \"33\u00d768=2244\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Pattern not found: $oldText\"\n    }\n}\n"}
